$wb = $excel.ActiveWorkbook

# Write a literal text value into a cell without Excel re-interpreting a
# numeric-looking string as a number (and without leaving a stray
# quotePrefix style behind, as a leading "'" on .Value would): build the
# text via a `="..."` formula, then copy/paste-special as values only so
# the formula collapses to a plain literal string.
function Set-TextValue($range, [string]$text) {
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# NOTE: the workbook has two sheets whose names differ only by case
# ("Vector_bf" vs "Vector_BF"); Worksheets.Item(<name>) resolves
# case-insensitively in this host, so sheets are addressed by their
# (1-based) tab index instead to avoid ambiguity.
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

$wsFollower = $wb.Worksheets.Item(3)
$wsPoint    = $wb.Worksheets.Item(4)
$wsVecbf    = $wb.Worksheets.Item(5)
$wsVecBF    = $wb.Worksheets.Item(6)

# Restricciones_del_follower (re-generated quadratic/linear problem values)
Set-TextValue $wsFollower.Range("A2") "-0.25 - x + y_1 + y_2"
Set-TextValue $wsFollower.Range("B2") "0.25"
Set-TextValue $wsFollower.Range("D2") "0.34"
Set-TextValue $wsFollower.Range("E2") "2.4"
Set-TextValue $wsFollower.Range("F2") "0.7000000000000001"

Set-TextValue $wsFollower.Range("A3") "1.7000000000000002 - y_1"
Set-TextValue $wsFollower.Range("B3") "-1.7000000000000002"
Set-TextValue $wsFollower.Range("D3") "0.14"
Set-TextValue $wsFollower.Range("E3") "8.0"
Set-TextValue $wsFollower.Range("F3") "2.1"

Set-TextValue $wsFollower.Range("A4") "-3.8 - y_2"
Set-TextValue $wsFollower.Range("B4") "-3.8"
Set-TextValue $wsFollower.Range("D4") "0.38"
Set-TextValue $wsFollower.Range("E4") "0.5"
Set-TextValue $wsFollower.Range("F4") "3.5"

# Punto_modificado
Set-TextValue $wsPoint.Range("A2") "5.25"
Set-TextValue $wsPoint.Range("B2") "1.7000000000000002"
Set-TextValue $wsPoint.Range("C2") "3.8"

# Vector_bf
Set-TextValue $wsVecbf.Range("A2") "-1.9000000000000001"
Set-TextValue $wsVecbf.Range("A3") "-0.96"

# Vector_BF
Set-TextValue $wsVecBF.Range("A2") "2.4"
Set-TextValue $wsVecBF.Range("A3") "5.8999999999999995"
Set-TextValue $wsVecBF.Range("A4") "-3.6999999999999997"

$excel.CutCopyMode = 0

Write-Host "edit applied"
